$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 591, shifting existing rows 591..637 down to 592..638.
$ws.Rows.Item(591).Insert()

# Populate the newly inserted row 591 with the new weekly record.
$ws.Cells.Item(591, 1).Value  = 9
$ws.Cells.Item(591, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(591, 3).Value  = "Metropolitana"
$ws.Cells.Item(591, 4).Value  = 45021
$ws.Cells.Item(591, 5).Value  = 13
$ws.Cells.Item(591, 6).Value  = "Fruta"
$ws.Cells.Item(591, 7).Value  = 100108
$ws.Cells.Item(591, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(591, 9).Value  = 100108002
$ws.Cells.Item(591, 10).Value = "Mango"
$ws.Cells.Item(591, 11).Value = "Sin especificar"
$ws.Cells.Item(591, 12).Value = "Primera"
$ws.Cells.Item(591, 13).Value = 660
$ws.Cells.Item(591, 14).Value = 7000
$ws.Cells.Item(591, 15).Value = 7500
$ws.Cells.Item(591, 16).Value = 7212
$ws.Cells.Item(591, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(591, 18).Value = "Perú"
$ws.Cells.Item(591, 19).Value = 1803
$ws.Cells.Item(591, 20).Value = 4
